$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'26.744.36"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +0.45%  "

$ws.Cells.Item(3, 4).Value = "'1.642.53"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -0.07%  "

$ws.Cells.Item(4, 5).Value = "  +0.39%  "

$ws.Cells.Item(5, 4).Value = "'216.51"
$ws.Cells.Item(5, 4).Style = "Normal"

$ws.Cells.Item(6, 4).Value = "'0.499"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.75%  "

$ws.Cells.Item(7, 5).Value = "  +0.27%  "

$ws.Cells.Item(8, 4).Value = "'0.0629"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +0.50%  "

$ws.Cells.Item(9, 5).Value = "  -0.43%  "

$ws.Cells.Item(10, 4).Value = "'19.13"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -0.50%  "

$ws.Cells.Item(11, 5).Value = "  -0.25%  "

$ws.Cells.Item(12, 4).Value = "'1.866.21"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -0.36%  "

$ws.Cells.Item(13, 4).Value = "'1.653.35"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +0.68%  "

$ws.Cells.Item(14, 5).Value = "  -1.52%  "

$ws.Cells.Item(15, 5).Value = "  -1.07%  "

$ws.Cells.Item(16, 4).Value = "'64.36"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -2.39%  "

$ws.Cells.Item(17, 4).Value = "'26.732.79"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +0.21%  "

$ws.Cells.Item(18, 5).Value = "  -1.78%  "

$ws.Cells.Item(19, 4).Value = "'213.53"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -2.27%  "

$ws.Cells.Item(20, 5).Value = "  +0.27%  "

$ws.Cells.Item(21, 4).Value = "'4.36"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -0.23%  "

$ws.Cells.Item(22, 4).Value = "'2.43"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +12.87%  "

$ws.Cells.Item(23, 4).Value = "'6.25"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.95%  "

$ws.Cells.Item(24, 5).Value = "  -2.36%  "

$ws.Cells.Item(25, 4).Value = "'144.95"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.94%  "

$ws.Cells.Item(26, 5).Value = "  +0.37%  "

$ws.Cells.Item(27, 4).Value = "'0.119"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -1.41%  "

$ws.Cells.Item(28, 5).Value = "  -0.33%  "

$ws.Cells.Item(29, 4).Value = "'15.63"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -1.38%  "

$ws.Cells.Item(30, 5).Value = "  -1.49%  "

$ws.Cells.Item(31, 5).Value = "  +0.34%  "

$ws.Cells.Item(32, 4).Value = "'3.31"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -2.38%  "

$ws.Cells.Item(33, 5).Value = "  -2.29%  "

$ws.Cells.Item(34, 4).Value = "'1.293.93"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +1.24%  "

$ws.Cells.Item(35, 4).Value = "'1.53"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -0.53%  "

$ws.Cells.Item(36, 5).Value = "  +1.37%  "

$ws.Cells.Item(37, 4).Value = "'0.0174"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -4.76%  "

$ws.Cells.Item(38, 4).Value = "'0.532"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +1.03%  "

$ws.Cells.Item(39, 5).Value = "  -0.35%  "

$ws.Cells.Item(40, 5).Value = "  +0.26%  "

$ws.Cells.Item(41, 4).Value = "'0.807"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -0.04%  "

$ws.Cells.Item(42, 5).Value = "  -0.16%  "

$ws.Cells.Item(43, 5).Value = "  -2.03%  "

$ws.Cells.Item(44, 4).Value = "'1.793.42"
$ws.Cells.Item(44, 4).Style = "Normal"

$ws.Cells.Item(45, 4).Value = "'61.60"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +3.00%  "

$ws.Cells.Item(46, 4).Value = "'91.35"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -1.85%  "

$ws.Cells.Item(47, 5).Value = "  -0.50%  "

$ws.Cells.Item(48, 4).Value = "'0.0522"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +0.99%  "

$ws.Cells.Item(49, 4).Value = "'7.66"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -1.40%  "

$ws.Cells.Item(50, 4).Value = "'0.0974"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.37%  "

$ws.Cells.Item(51, 2).Value = "Mantle"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(51, 4).Value = "'0.407"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.09%  "
